# Readme-reg.xlsx update: rename scenario sheets, add a footnote cell,
# and move the active selection to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the six scenario sheets -------------------------------------
$wb.Worksheets.Item(1).Name = "시나리오 1. 생산량 예측 (단순 선형)"
$wb.Worksheets.Item(2).Name = "시나리오 1.1. 생산량 예측 (Auto, Neural)"
$wb.Worksheets.Item(3).Name = "시나리오 2. 유사제품수요증감예측 "
$wb.Worksheets.Item(4).Name = "시나리오 2.1 유사제품월별수요증감예측 "
$wb.Worksheets.Item(5).Name = "시나리오 3. 월평균 예측 (선형분석)"
$wb.Worksheets.Item(6).Name = "시나리오 4.  유사도 분석 (거리계산)"

# --- 2. Note that the simple-linear model output is only produced for the
#        simple-linear scenario ---------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("H15").Value = "(단순 선형에서만 산출)"

# --- 3. Restore each sheet's scroll position / selection -------------------
$ws1.Range("F23").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F19").Select()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("D30").Select()

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").Select()

# --- 4. Make the last scenario sheet the active tab -------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Activate()
$ws6.Range("A3").Select()
